$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1157
$ws1.Range("F7").Value = 12216
$ws1.Range("F11").Value = 12011
$ws1.Range("F13").Value = 2628
$ws1.Range("F15").Value = 45
$ws1.Range("F19").Value = 358

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1157
$ws4.Range("F9").Value = 12216
$ws4.Range("F13").Value = 12011
$ws4.Range("F15").Value = 2628
$ws4.Range("F17").Value = 45
$ws4.Range("F21").Value = 358
